# FotogroepWaalre.xlsx - add two new rows of source-code-line-count
# tracking data (2023-02-19 and 2023-02-20) to the "Table" sheet.
#
# Row 149: date 44976 (2023-02-19)
# Row 150: date 44977 (2023-02-20)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

# Reference cell (last existing data row) used to copy over its
# (column-driven) number styling for the two computed-total cells,
# which otherwise pick up a stray auto-applied number format.
$refTotalStyle = $ws.Cells.Item(148, 8).Style

# ---- Row 149 -----------------------------------------------------
$ws.Cells.Item(149, 1).Value = 44976
$ws.Cells.Item(149, 2).Value = 9
$ws.Cells.Item(149, 3).Formula = "=SUM(D149:F149)"
$ws.Cells.Item(149, 4).Value = 1105
$ws.Cells.Item(149, 5).Value = 371
$ws.Cells.Item(149, 6).Value = 197
$ws.Cells.Item(149, 7).Value = 52
$ws.Cells.Item(149, 8).Formula = "=SUM(I149:K149)"
$ws.Cells.Item(149, 8).Style = $refTotalStyle
$ws.Cells.Item(149, 9).Value = 3724
$ws.Cells.Item(149, 10).Value = 544
$ws.Cells.Item(149, 11).Value = 476

# ---- Row 150 -----------------------------------------------------
$ws.Cells.Item(150, 1).Value = 44977
$ws.Cells.Item(150, 2).Value = 9
$ws.Cells.Item(150, 3).Formula = "=SUM(D150:F150)"
$ws.Cells.Item(150, 4).Value = 1105
$ws.Cells.Item(150, 5).Value = 371
$ws.Cells.Item(150, 6).Value = 197
$ws.Cells.Item(150, 7).Value = 52
$ws.Cells.Item(150, 8).Formula = "=SUM(I150:K150)"
$ws.Cells.Item(150, 8).Style = $refTotalStyle
$ws.Cells.Item(150, 9).Value = 4747
$ws.Cells.Item(150, 10).Value = 544
$ws.Cells.Item(150, 11).Value = 476

# Park the visible selection on the Table sheet at I151 (just past the
# newly-added data) without leaving "Table" as the active sheet -
# the workbook was (and stays) opened on "Graph".
$wsGraph = $wb.Worksheets.Item("Graph")
[void]$ws.Activate()
[void]$ws.Range("I151").Select()
[void]$wsGraph.Activate()
